$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '55.359.45'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -5.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.889.40'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -5.72%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '485.67'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -7.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.09'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -8.02%  '
$ws.Range('E8').Value = '  -7.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.09'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -5.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.344'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -6.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.381.07'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -5.72%  '
$ws.Range('E13').Value = '  -4.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.52'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -5.02%  '
$ws.Range('E15').Value = '  -9.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '55.359.18'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -4.94%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.92'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -4.91%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.885.36'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -6.11%  '
$ws.Range('E19').Value = '  -6.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.56'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -7.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '311.24'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -8.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').ClearFormats()
$ws.Range('E24').Value = '  -5.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '61.70'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -5.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.997'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('E27').Value = '  -6.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0828'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -14.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.31'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -9.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.90'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -8.73%  '
$ws.Range('E31').Value = '  -6.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.51'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -7.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.10'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -10.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '148.81'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -5.48%  '
$ws.Range('E35').Value = '  -9.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.54'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -7.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '24.19'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.17'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -10.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0645'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -7.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.920.26'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -5.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.04'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.63%  '
$ws.Range('E43').Value = '  -8.15%  '
$ws.Range('E44').Value = '  -6.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.080.95'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -10.95%  '
$ws.Range('E46').Value = '  -10.31%  '
$ws.Range('E47').Value = '  -5.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.901'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -11.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0226'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -6.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.45'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -7.21%  '
$ws.Range('E51').Value = '  -7.80%  '
